$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabla Pivote")

# Insert the 7 new "_Incremento" columns, processed right-to-left so the
# column letter used for each EntireColumn.Insert() call (computed against
# the ORIGINAL 8-column layout: A Periodo, B Aula, C Lab.Computadoras,
# D Lab.Fisica, E Lab.Quimica, F Taller, G Virtual, H Total) stays valid
# while earlier (left-hand) columns are still untouched.
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("H1").EntireColumn.Insert()
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("C1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("C1").Value = "Aula_Incremento"
$ws.Range("E1").Value = "Laboratorio de Computadoras_Incremento"
$ws.Range("G1").Value = "Laboratorio de Física_Incremento"
$ws.Range("I1").Value = "Laboratorio de Química_Incremento"
$ws.Range("K1").Value = "Taller_Incremento"
$ws.Range("M1").Value = "Virtual_Incremento"
$ws.Range("O1").Value = "Total_Incremento"

# Values for the new columns, one row per period (rows 2-21), in the
# column order C, E, G, I, K, M, O.
$incrementoData = @(
    @(4, 0, 0, 0, 4, 0, 8),
    @(22, 4, 0, 0, 0, 3, 29),
    @(14, 0, 0, 3, 4, 5, 26),
    @(21, 0, 0, 2, 0, 5, 28),
    @(18, 0, 0, 2, 0, 0, 20),
    @(32, 3, 5, 3, 0, 3, 46),
    @(26, 0, 3, 0, 0, 0, 29),
    @(12, 0, 0, 3, 4, 0, 19),
    @(32, 0, 0, 0, 0, 4, 36),
    @(40, 0, 0, 0, 0, 4, 44),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0)
)

$newCols = @(3, 5, 7, 9, 11, 13, 15)

for ($i = 0; $i -lt $incrementoData.Count; $i++) {
    $rowValues = $incrementoData[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $newCols.Count; $j++) {
        $ws.Cells.Item($rowNum, $newCols[$j]).Value = $rowValues[$j]
    }
}
